$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 11: replaces old one-sample t-test row with new QQ-plot question ---
# Row 11
$ws.Range("A11").Value = 'Reasonable QQ plots'
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 'which qqplots seem normally distributed'
$ws.Range("D11").Value = 'lm-normality-plots'
$ws.Range("E11").Value = ''

# Row 12
$ws.Range("A12").Value = ''
$ws.Range("B12").Value = ''
$ws.Range("C12").Value = ''
$ws.Range("D12").Value = ''
$ws.Range("E12").Value = ''

# Row 13
$ws.Range("A13").Value = ''
$ws.Range("B13").Value = ''
$ws.Range("C13").Value = ''
$ws.Range("D13").Value = ''
$ws.Range("E13").Value = ''

# Row 14
$ws.Range("A14").Value = ''
$ws.Range("B14").Value = ''
$ws.Range("C14").Value = ''
$ws.Range("D14").Value = ''
$ws.Range("E14").Value = ''

# Row 15
$ws.Range("A15").Value = ''
$ws.Range("B15").Value = ''
$ws.Range("C15").Value = ''
$ws.Range("D15").Value = ''
$ws.Range("E15").Value = ''

# Row 16
$ws.Range("A16").Value = 'one-sample t-test'
$ws.Range("B16").Value = 5
$ws.Range("C16").Value = 'read test-statistic, read sided t-test'
$ws.Range("D16").Value = 'schoice-2samtt-interpret-t, schoice-2samtt-interpret-alt, schoice-2samtt-interpret-p'
$ws.Range("E16").Value = 'single choice'

# Row 17
$ws.Range("A17").Value = 'assumptions of one-sample t-test'
$ws.Range("B17").Value = 5
$ws.Range("C17").Value = ''
$ws.Range("D17").Value = ''
$ws.Range("E17").Value = ''

# Row 18
$ws.Range("A18").Value = 'two-sample t-test'
$ws.Range("B18").Value = 5
$ws.Range("C18").Value = 'know which scenario to use t-test'
$ws.Range("D18").Value = 'which-ttest-to-use'
$ws.Range("E18").Value = 'schoice'

# Row 19
$ws.Range("A19").Value = 'assumptions of two-sample t-test'
$ws.Range("B19").Value = 5
$ws.Range("C19:E19").Clear()

# Row 20
$ws.Range("A20").Value = 'matched pairs t-test'
$ws.Range("B20").Value = 5
$ws.Range("C20").Value = ''
$ws.Range("D20").Value = ''
$ws.Range("E20").Value = ''

# Row 21
$ws.Range("A21").Value = 'assumptions of matched-pairs t-test'
$ws.Range("B21").Value = 5
$ws.Range("C21").Value = ''
$ws.Range("D21").Value = ''
$ws.Range("E21").Value = ''

# Row 22
$ws.Range("A22").Value = 'chi-squared test'
$ws.Range("B22").Value = 6
$ws.Range("C22").Value = ''
$ws.Range("D22").Value = ''
$ws.Range("E22").Value = ''

# Row 23
$ws.Range("A23").Value = 'interpret chi-squared test'
$ws.Range("B23").Value = 6
$ws.Range("C23").Value = ''
$ws.Range("D23").Value = ''
$ws.Range("E23").Value = ''

# Row 24
$ws.Range("A24").Value = 'mann-whitney test'
$ws.Range("B24").Value = 6
$ws.Range("C24").Value = 'interpret slope, interpret intercept, assumptions of lm check'
$ws.Range("D24").Value = 'lm-slope-properties, lm-intercept-properties, lm-assmptions'
$ws.Range("E24").Value = 'schoice'

# Row 25
$ws.Range("A25").Value = 'interpret mann-whiteny test'
$ws.Range("B25").Value = 6
$ws.Range("C25").Value = ''
$ws.Range("D25").Value = ''
$ws.Range("E25").Value = ''

# Row 26
$ws.Range("A26").Value = 'linear regression (theoretical) '
$ws.Range("B26").Value = 7
$ws.Range("C26").Value = 'slope estimate'
$ws.Range("D26").Value = 'linear-model-slope'
$ws.Range("E26").Value = 'num'

# Row 27
$ws.Range("A27").Value = 'produce diagnostic residual plots in r'
$ws.Range("B27").Value = 7
$ws.Range("C27").Value = 'which plot to assess? & which seems reasonable? : linearity, normality'
$ws.Range("D27").Value = 'lm-linearity-check, lm-linearity-plots, lm-normality-check, lm-normality-plots, lm-homoscedasticity-check, lm-homoscedasticity-plots'
$ws.Range("E27").Value = 'schoice'

# Row 28
$ws.Range("A28").Value = 'interpret coefficient estimates in lm'
$ws.Range("B28").Value = 7
$ws.Range("C28").Value = ''
$ws.Range("D28").Value = ''
$ws.Range("E28").Value = ''

# Row 29
$ws.Range("A29").Value = 'assess regression assumptions for lm'
$ws.Range("B29").Value = 7
$ws.Range("C29").Value = ''
$ws.Range("D29").Value = ''
$ws.Range("E29").Value = ''

# Row 30
$ws.Range("A30").Value = 'interpret CI & PI from a simple lm model'
$ws.Range("B30").Value = 7
$ws.Range("C30").Value = ''
$ws.Range("D30").Value = ''
$ws.Range("E30").Value = ''

# --- New blank placeholder rows 31-32 (full A:E) and 33-34 (A:B only) ---
$ws.Range("A31:E32").WrapText = $true
$ws.Range("A33:B34").WrapText = $true

# --- Row heights for rows 26-30 (auto-computed wrap height in Excel; set explicitly here) ---
$ws.Rows.Item(26).RowHeight = 17
$ws.Rows.Item(27).RowHeight = 68
$ws.Rows.Item(28).RowHeight = 17
$ws.Rows.Item(29).RowHeight = 17
$ws.Rows.Item(30).RowHeight = 17

# --- Selection / view state ---
$ws.Range("D12").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1